$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# "soundkonverter added for Manjaro"
#
# A new row is inserted right after "soundconverter" (row 32) for the new
# program "soundkonverter": supported on Manjaro only (Windows / Debian not
# supported). This pushes every following program row down by one.
#
# A second new row is appended at the very end of the table for "wamp"
# (supported on Windows only).
# ---------------------------------------------------------------------------

# Reference cells that already carry the two fill styles used throughout the
# table, so the new cells get the exact same formatting as the rest of the
# sheet (green = supported, red = not supported).
$supported = $ws.Range("B4")       # green "supported" style
$notSupported = $ws.Range("B5")    # red "not supported" style

# 1) Insert a new row at 33 (shifts old rows 33-50 down to 34-51) and fill it
#    in with the soundkonverter data.
$ws.Rows("33:33").Insert()

$ws.Range("A33").Value = "soundkonverter"
$notSupported.Copy($ws.Range("B33"))   # Windows -> not supported
$supported.Copy($ws.Range("C33"))      # Manjaro -> supported
$notSupported.Copy($ws.Range("D33"))   # Debian  -> not supported

# 2) Append a new row (52) for wamp at the bottom of the table.
$ws.Range("A52").Value = "wamp"
$supported.Copy($ws.Range("B52"))      # Windows -> supported
$notSupported.Copy($ws.Range("C52"))   # Manjaro -> not supported
$notSupported.Copy($ws.Range("D52"))   # Debian  -> not supported

# Keep the view roughly where the original edit left it.
$ws.Range("C50").Select()
$excel.ActiveWindow.ScrollRow = 28

Write-Output "done"
